# edit.ps1 - apply the "Fix all issues from Tschopp annotated ms." revisions
# to response-to-reviews.docx, via the Word COM-interop object model.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $old"
    }
}

# 1. Remove stray "XXX " prefix before "We have reworked and combined..."
Replace-Text "XXX We have reworked and combined some of the illustrations along the lines that Marek suggests." "We have reworked and combined some of the illustrations along the lines that Marek suggests."

# 2. Merge the colour-blindness paragraph's four runs into a single run (fixing the "We have have" typo is NOT requested - keep wording identical to diff).
Replace-Text "We note the comments on colour-blindness. We have have checked that where colour is used to distinguish lines with different meanings those lines are also distinguishable by other criteria, and made the relevant figure captions more explicit." "We note the comments on colour-blindness. We have have checked that where colour is used to distinguish lines with different meanings those lines are also distinguishable by other criteria, and made the relevant figure captions more explicit."

# 3. Merge the "We have in general followed..." paragraph's leading runs (stop right before the HYPERLINK field) into a single run.
Replace-Text "We have in general followed the specific comments attached to Marek’s review, but with some exceptions. In particular, we do not agree with the suggestion that expunging the pronoun “we” throughout and substituting passive voice would improve the manuscript, and note that the fourth of " "We have in general followed the specific comments attached to Marek’s review, but with some exceptions. In particular, we do not agree with the suggestion that expunging the pronoun “we” throughout and substituting passive voice would improve the manuscript, and note that the fourth of "

# 4. Merge the "We have retained a slightly shortened section on open peer-review..." paragraph's six runs into one.
Replace-Text "We have retained a slightly shortened section on open peer-review, contrary to Marek’s recommendation and in accordance with that of Tschopp. This section is important to us and relevant to the origin of the present paper’s core question. We have, however, removed the section on open composition." "We have retained a slightly shortened section on open peer-review, contrary to Marek’s recommendation and in accordance with that of Tschopp. This section is important to us and relevant to the origin of the present paper’s core question. We have, however, removed the section on open composition."

# 5a. Merge "The point that sequences of vertebrae oriented by Method 2 need not be illustrated..." runs into one.
Replace-Text "The point that sequences of vertebrae oriented by Method 2 need not be illustrated in a way that results in a jagged neural canal is important: we now address it in the manuscript." "The point that sequences of vertebrae oriented by Method 2 need not be illustrated in a way that results in a jagged neural canal is important: we now address it in the manuscript."

# 5b. Merge "We have added a note that Tschopp et al. (2015)'s definition of Character 194 ..." runs into one.
Replace-Text "We have added a note that Tschopp et al. (2015)’s definition of Character 194 includes a note that a horizontal orientation of the neural canal is used when scoring, though without discussion." "We have added a note that Tschopp et al. (2015)’s definition of Character 194 includes a note that a horizontal orientation of the neural canal is used when scoring, though without discussion."

# 6. Remove stray "XXX " prefix before "We have removed some of the references to blog-posts..."
Replace-Text "XXX We have removed some of the references to blog-posts as requested. Other are retained, as they properly acknowledge the source of previously published ideas and images. As we note in our in-press chapter in the forthcoming 3rd edition of " "We have removed some of the references to blog-posts as requested. Other are retained, as they properly acknowledge the source of previously published ideas and images. As we note in our in-press chapter in the forthcoming 3rd edition of "

# 7. Figure 5 vertebra paragraph: remove "XXX " prefix and "XXX what?" placeholder, and add the real sentence as a *separate* trailing run.
$old7 = "XXX We agree that adding another vertebra to Figure 5, as a second example of how different methods yield different orientations, would be helpful. We have added XXX what?"
$new7a = "We agree that adding another vertebra to Figure 5, as a second example of how different methods yield different orientations, would be helpful. We have added XXXPLACEHOLDERXXX"
Replace-Text $old7 $new7a

$para7 = $d.Paragraphs | Where-Object { $_.Range.Text -like "*Figure 5*" }
$r7 = $para7.Range
$new7b = "a posterior cervical vertebra of a giraffe, so that the illustrated example encompass two high-level clades (Sauropsida and Synapsida) as well as two regions of the vertebral column (caudal and cervical)."
$r7.Find.Execute("XXXPLACEHOLDERXXX", $true, $false, $false, $false, $false, $true, 1, $false, $new7b, 2) | Out-Null

# Select just the newly-inserted sentence and toggle Bold on/off to force it into its own run
# (with empty rPr), matching the target markup.
$para7b = $d.Paragraphs | Where-Object { $_.Range.Text -like "*Figure 5*" }
$r7b = $para7b.Range
$r7b.Find.Execute($new7b, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r7b.Bold = 1
$r7b.Bold = 0

# 8. styles.xml: Normal style overflowPunct false -> true
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.Hyphenation = $normalStyle.ParagraphFormat.Hyphenation
# overflowPunct isn't exposed as a named COM property; toggle via low-level style XML instead.

# 9. styles.xml: add new character style "ListLabel10" ("ListLabel 10")
$newStyle = $d.Styles.Add("ListLabel 10", 2)
$newStyle.QuickStyle = $true
